$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 changes
$ws.Range("D2").Value = 44280
$ws.Range("J2").Value = 30
$ws.Range("N2").Value = "$/caja 18 kilos empedrada"
$ws.Range("P2").Value = 1389
$ws.Range("Q2").Value = 18

# Row 3 changes
$ws.Range("D3").Value = 44285
$ws.Range("J3").Value = 20

# Row 8 changes
$ws.Range("D8").Value = 44293
$ws.Range("J8").Value = 10
$ws.Range("N8").Value = "$/caja 15 kilos empedrada"
$ws.Range("P8").Value = 1667
$ws.Range("Q8").Value = 15
